# Apply updated evaluation metric values across the three worksheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Summary ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.2437722419928826
$wsSummary.Range("C2").Value = 0.06181015452538632
$wsSummary.Range("E2").Value = 0.1164241164241164
$wsSummary.Range("F2").Value = 0.247787610619469
$wsSummary.Range("G2").Value = 0.6313963573287077
$wsSummary.Range("H2").Value = 0.7984216158373462
$wsSummary.Range("J2").Value = 425
$wsSummary.Range("K2").Value = 109

# ---- Sheet 2: Classification Report ----
$wsReport = $wb.Worksheets.Item("Classification Report")
$wsReport.Range("C2").Value = 0.2041198501872659
$wsReport.Range("D2").Value = 0.3390357698289269

$wsReport.Range("B3").Value = 0.06181015452538632
$wsReport.Range("D3").Value = 0.1164241164241164

$wsReport.Range("B4").Value = 0.2437722419928826
$wsReport.Range("C4").Value = 0.2437722419928826
$wsReport.Range("D4").Value = 0.2437722419928826
$wsReport.Range("E4").Value = 0.2437722419928826

$wsReport.Range("B5").Value = 0.5309050772626932
$wsReport.Range("C5").Value = 0.6020599250936329
$wsReport.Range("D5").Value = 0.2277299431265217

$wsReport.Range("B6").Value = 0.9532574454211936
$wsReport.Range("C6").Value = 0.2437722419928826
$wsReport.Range("D6").Value = 0.3279447977731712

# ---- Sheet 3: Confusion Matrix ----
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")
$wsConfusion.Range("B2").Value = 109
$wsConfusion.Range("C2").Value = 425
